{"js": "// R\u00e9sum\u00e9 (\u540e\u7aef\u7b80\u5386.docx) edits:\n//  1) \"\u6846\u67b6\" bullet: \"\uff0c\u719f\u6089 Restful\u3002 \" -> \"\uff0c\u719f\u6089Gin\uff0cGorm\u3002\"\n//     (drop the trailing \"Restful\" skill in favor of \"Gin\uff0cGorm\")\n//  2) Awards bullet: \"CET4\" -> \"CET6\" (keep the long run of trailing\n//     spaces / tabs intact) and the trailing date \"2023-08-24\" -> \"393\"\n\nconst body = context.document.body;\n\n// --- 1) Replace \" Restful\u3002 \" (and the \"\uff0c\u719f\u6089\" immediately preceding it\n//        stays untouched) with \"Gin\uff0cGorm\u3002\" -------------------------------\nconst restfulHits = body.search(\"\uff0c\u719f\u6089 Restful\u3002 \", { matchCase: true });\nrestfulHits.load(\"text\");\nawait context.sync();\n\nif (restfulHits.items.length > 0) {\n  restfulHits.items[0].insertText(\"\uff0c\u719f\u6089Gin\uff0cGorm\u3002\", Word.InsertLocation.replace);\n}\n\n// --- 2) \"CET4\" -> \"CET6\" ----------------------------------------------------\n// Find the unique \"CET4\" span, split it so we isolate just the \"4\"\n// character (leaving \"CET\" and the long trailing space run alone), then\n// swap that single character for \"6\".\nconst cetHits = body.search(\"CET4\", { matchCase: true });\ncetHits.load(\"text\");\nawait context.sync();\n\nif (cetHits.items.length > 0) {\n  const cetHit = cetHits.items[0];\n  const parts = cetHit.split([\"T\"], false, false);\n  parts.load(\"text\");\n  await context.sync();\n  // parts.items -> [\"CET\", \"4\"]\n  const digitRange = parts.items[parts.items.length - 1];\n  digitRange.insertText(\"6\", Word.InsertLocation.replace);\n}\n\n// --- 3) \"2023-08-24\" -> \"393\" ----------------------------------------------\nconst dateHits = body.search(\"2023-08-24\", { matchCase: true });\ndateHits.load(\"text\");\nawait context.sync();\n\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"393\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# R\u00e9sum\u00e9 (\u540e\u7aef\u7b80\u5386.docx) edits:\n#  1) \"\u6846\u67b6\" bullet: \"\uff0c\u719f\u6089 Restful\u3002 \" -> \"\uff0c\u719f\u6089Gin\uff0cGorm\u3002\"\n#     (drop the trailing \"Restful\" skill in favor of \"Gin\uff0cGorm\")\n#  2) Awards bullet: \"CET4\" -> \"CET6\" (leave the long run of trailing\n#     spaces / tabs intact) and the trailing date \"2023-08-24\" -> \"393\"\n\n$d = $word.ActiveDocument\n\n# --- 1) \"\uff0c\u719f\u6089 Restful\u3002 \" -> \"\uff0c\u719f\u6089Gin\uff0cGorm\u3002\" --------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"\uff0c\u719f\u6089 Restful\u3002 \"\n$rng.Find.Replacement.Text = \"\uff0c\u719f\u6089Gin\uff0cGorm\u3002\"\n$rng.Find.Execute($rng.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng.Find.Replacement.Text, 2) | Out-Null\n\n# --- 2) \"CET4\" -> \"CET6\" -----------------------------------------------------\n# Locate the unique \"CET4\" span, then isolate just the trailing \"4\"\n# character (leaving \"CET\" and the long run of trailing spaces alone) and\n# swap that single character for \"6\".\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"CET4\"\n$found = $rng2.Find.Execute()\nif ($found) {\n  $digitRange = $d.Range($rng2.End - 1, $rng2.End)\n  $digitRange.Text = \"6\"\n}\n\n# --- 3) \"2023-08-24\" -> \"393\" ------------------------------------------------\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$rng3.Find.Text = \"2023-08-24\"\n$rng3.Find.Replacement.Text = \"393\"\n$rng3.Find.Execute($rng3.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng3.Find.Replacement.Text, 2) | Out-Null\n"}
